# Insert a new row above row 191. Excel shifts the existing rows 191-286
# down to 192-287, preserving all of their data (dimension grows to T287).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(191).Insert()

# The new (blank) row 191 should start life as a duplicate of the record
# that is now sitting in row 192 (the former row 191), then have the
# fields that differ for the new record overwritten below.
$ws.Range("A192:T192").Copy()
$ws.Range("A191").PasteSpecial()

# Apply the field values that are specific to the newly inserted record.
$ws.Range("D191").Value = 44992
$ws.Range("L191").Value = "Segunda"
$ws.Range("M191").Value = 50
$ws.Range("N191").Value = 25000
$ws.Range("O191").Value = 25000
$ws.Range("P191").Value = 25000
$ws.Range("Q191").Value = "$/caja 14 unidades"
$ws.Range("S191").Value = 1786
$ws.Range("T191").Value = 14
